# Daily attendance processing - 2025-11-14 04:50:49
# Normalize the "Recorded By" (column G) values: when the list of
# recorders starts with "System, ", move the "System" entry to the end
# of the comma-separated list instead of the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2
    if ($value -ne $null -and $value -is [string] -and $value.StartsWith("System, ")) {
        $rest = $value.Substring(8)
        $cell.Value = $rest + ", System"
    }
}
